$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-6 per repulled data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 1
